# Insert a new record at row 600 (Macroferia Regional de Talca - Brocoli),
# which pushes the existing rows 600:650 down to 601:651.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(600).Insert()

$ws.Range("A600").Value = 5
$ws.Range("B600").Value = "Macroferia Regional de Talca"
$ws.Range("C600").Value = "Maule"
$ws.Range("D600").Value = 45223
$ws.Range("E600").Value = 7
$ws.Range("F600").Value = 100112023
$ws.Range("G600").Value = "Brócoli"
$ws.Range("H600").Value = "Sin especificar"
$ws.Range("I600").Value = "Primera"
$ws.Range("J600").Value = 4000
$ws.Range("K600").Value = 900
$ws.Range("L600").Value = 900
$ws.Range("M600").Value = 900
$ws.Range("N600").Value = "$/unidad"
$ws.Range("O600").Value = "Región del Maule"
$ws.Range("P600").Value = 900
$ws.Range("Q600").Value = 1
$ws.Range("R600").Value = "Hortaliza"
